# "Natmi following Dr Hou advice"
#
# The NATMI ligand-receptor output for Lrfn3-Lrfn3 is regenerated with a
# third cluster ("ECs") added alongside the existing "FAPs" and "sCs"
# clusters. The 2x2 sending/target-cluster matrix (4 data rows) becomes a
# full 3x3 matrix (9 data rows), and every derived statistic is recomputed
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (A:T), one row per Sending/Target
# cluster combination of {ECs, FAPs, sCs}. Ligand symbol (B) and Receptor
# symbol (C) are always "Lrfn3" for this pair.
$rows = @(
  @{ row=2;  A="ECs";  B="Lrfn3"; C="Lrfn3"; D="ECs";  E=1; F=0.3333333333333333; G=0.04467366666666667; H=0.134021;           I=0.009984424503391825; J=0.009984424503391825; K=1; L=0.3333333333333333; M=0.04467366666666667; N=0.134021;           O=0.009984424503391825; P=0.009984424503391825; Q=0.001995736493444444; R=0.017961628441;    S=0.0000996887326639311; T=0.0000996887326639311 },
  @{ row=3;  A="ECs";  B="Lrfn3"; C="Lrfn3"; D="FAPs"; E=1; F=0.3333333333333333; G=0.04467366666666667; H=0.134021;           I=0.009984424503391825; J=0.009984424503391825; K=3; L=1;                  M=2.659496333333333;  N=7.978489;            O=0.5943890962732866;   P=0.5943890962732866;   Q=0.1188094526965556;   R=1.069285074269;     S=0.005934633057379925;  T=0.005934633057379925  },
  @{ row=4;  A="ECs";  B="Lrfn3"; C="Lrfn3"; D="sCs";  E=1; F=0.3333333333333333; G=0.04467366666666667; H=0.134021;           I=0.009984424503391825; J=0.009984424503391825; K=3; L=1;                  M=1.770165666666667;  N=5.310497;            O=0.3956264792233216;   P=0.3956264792233216;   Q=0.07907979093744444;  R=0.7117181184369999; S=0.003950102713347969;  T=0.003950102713347969  },
  @{ row=5;  A="FAPs"; B="Lrfn3"; C="Lrfn3"; D="ECs";  E=3; F=1;                  G=2.659496333333333;  H=7.978489;            I=0.5943890962732866;   J=0.5943890962732866;   K=1; L=0.3333333333333333; M=0.04467366666666667; N=0.134021;           O=0.009984424503391825; P=0.009984424503391825; Q=0.1188094526965556;   R=1.069285074269;     S=0.005934633057379925;  T=0.005934633057379925  },
  @{ row=6;  A="FAPs"; B="Lrfn3"; C="Lrfn3"; D="FAPs"; E=3; F=1;                  G=2.659496333333333;  H=7.978489;            I=0.5943890962732866;   J=0.5943890962732866;   K=3; L=1;                  M=2.659496333333333;  N=7.978489;            O=0.5943890962732866;   P=0.5943890962732866;   Q=7.072920747013444;    R=63.65628672312099;  S=0.3532983977685744;    T=0.3532983977685744    },
  @{ row=7;  A="FAPs"; B="Lrfn3"; C="Lrfn3"; D="sCs";  E=3; F=1;                  G=2.659496333333333;  H=7.978489;            I=0.5943890962732866;   J=0.5943890962732866;   K=3; L=1;                  M=1.770165666666667;  N=5.310497;            O=0.3956264792233216;   P=0.3956264792233216;   Q=4.707749099892555;    R=42.369741899033;    S=0.2351560654473323;    T=0.2351560654473323    },
  @{ row=8;  A="sCs";  B="Lrfn3"; C="Lrfn3"; D="ECs";  E=3; F=1;                  G=1.770165666666667;  H=5.310497;            I=0.3956264792233216;   J=0.3956264792233216;   K=1; L=0.3333333333333333; M=0.04467366666666667; N=0.134021;           O=0.009984424503391825; P=0.009984424503391825; Q=0.07907979093744444;  R=0.7117181184369999; S=0.003950102713347969;  T=0.003950102713347969  },
  @{ row=9;  A="sCs";  B="Lrfn3"; C="Lrfn3"; D="FAPs"; E=3; F=1;                  G=1.770165666666667;  H=5.310497;            I=0.3956264792233216;   J=0.3956264792233216;   K=3; L=1;                  M=2.659496333333333;  N=7.978489;            O=0.5943890962732866;   P=0.5943890962732866;   Q=4.707749099892555;    R=42.369741899033;    S=0.2351560654473323;    T=0.2351560654473323    },
  @{ row=10; A="sCs";  B="Lrfn3"; C="Lrfn3"; D="sCs";  E=3; F=1;                  G=1.770165666666667;  H=5.310497;            I=0.3956264792233216;   J=0.3956264792233216;   K=3; L=1;                  M=1.770165666666667;  N=5.310497;            O=0.3956264792233216;   P=0.3956264792233216;   Q=3.133486487445444;    R=28.201378387009;    S=0.1565203110626413;    T=0.1565203110626413    }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.row, 1).Value  = $r.A
    $ws.Cells.Item($r.row, 2).Value  = $r.B
    $ws.Cells.Item($r.row, 3).Value  = $r.C
    $ws.Cells.Item($r.row, 4).Value  = $r.D
    $ws.Cells.Item($r.row, 5).Value  = $r.E
    $ws.Cells.Item($r.row, 6).Value  = $r.F
    $ws.Cells.Item($r.row, 7).Value  = $r.G
    $ws.Cells.Item($r.row, 8).Value  = $r.H
    $ws.Cells.Item($r.row, 9).Value  = $r.I
    $ws.Cells.Item($r.row, 10).Value = $r.J
    $ws.Cells.Item($r.row, 11).Value = $r.K
    $ws.Cells.Item($r.row, 12).Value = $r.L
    $ws.Cells.Item($r.row, 13).Value = $r.M
    $ws.Cells.Item($r.row, 14).Value = $r.N
    $ws.Cells.Item($r.row, 15).Value = $r.O
    $ws.Cells.Item($r.row, 16).Value = $r.P
    $ws.Cells.Item($r.row, 17).Value = $r.Q
    $ws.Cells.Item($r.row, 18).Value = $r.R
    $ws.Cells.Item($r.row, 19).Value = $r.S
    $ws.Cells.Item($r.row, 20).Value = $r.T
}
